{"js": "// Enhance Siege Analytics \"GIS & Geospatial Analysis Consulting\" entry with\n// three new bullet points describing boundary mapping, demographic miscoding\n// discovery, and boundary estimation tooling.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorText = \"GIS & Geospatial Analysis Consulting\";\nconst newBullets = [\n  \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n  \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n  \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n];\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// Insert the three bullets, one after another, immediately after the anchor\n// paragraph. Each new insert happens \"After\" the previous inserted paragraph\n// so the three bullets stay in order right below the anchor.\nlet insertAfter = anchorParagraph;\nfor (const bulletText of newBullets) {\n  insertAfter = insertAfter.insertParagraph(bulletText, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Enhance Siege Analytics \"GIS & Geospatial Analysis Consulting\" entry with\n# three new bullet points describing boundary mapping, demographic miscoding\n# discovery, and boundary estimation tooling.\n\n$d = $word.ActiveDocument\n\n$bullets = @(\n    \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n    \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n    \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n)\n\n$anchorText = \"GIS & Geospatial Analysis Consulting\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $anchorText) {\n        $anchor = $p\n        foreach ($bulletText in $bullets) {\n            $anchor.Range.InsertParagraphAfter()\n            $anchor = $anchor.Next()\n            $anchor.Range.Text = $bulletText\n        }\n        break\n    }\n}\n"}
